{"js": "// Replace the \"Some protected text.\" / \"End of demonstration.\" / trailing\n// empty paragraph block with a single bold, red error message appended to\n// the paragraph that hosts the `m:usercontent zone1` field.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Paragraph holding the <w:fldSimple w:instr=\"m:usercontent zone1\"/> field\n// (the first paragraph right after the intro sentence).\nconst fieldParagraph = paragraphs.items[1];\n\n// Append the new bold, red run right after the field, inside the same\n// paragraph.\nconst errorRange = fieldParagraph.insertText(\n  \"Invalid if statement: Unexpected tag EOF missing [ENDUSERDOC]\",\n  Word.InsertLocation.end\n);\nerrorRange.font.bold = true;\nerrorRange.font.color = \"#FF0000\";\n\n// Drop the next three paragraphs: \"Some protected text.\",\n// \"End of demonstration.\" and the trailing empty paragraph.\nparagraphs.items[2].delete();\nparagraphs.items[3].delete();\nparagraphs.items[4].delete();\n\nawait context.sync();\n", "ps1": "# Replace the \"Some protected text.\" / \"End of demonstration.\" / trailing\n# empty paragraph block with a single bold, red error message appended to\n# the paragraph that hosts the `m:usercontent zone1` field.\n\n$d = $word.ActiveDocument\n\n# Paragraph 2 holds <w:fldSimple w:instr=\"m:usercontent zone1\"/>; append the\n# new run right after the field, inside the same paragraph.\n$fieldParagraph = $d.Paragraphs.Item(2)\n$fieldParagraph.Range.InsertAfter(\"Invalid if statement: Unexpected tag EOF missing [ENDUSERDOC]\")\n\n# Locate the text we just inserted and make it bold + red.\n$errRange = $d.Content\n$errRange.Find.Execute(\"Invalid if statement: Unexpected tag EOF missing [ENDUSERDOC]\")\n$errRange.Font.Bold = 1\n$errRange.Font.Color = 255\n\n# Drop the next three paragraphs: \"Some protected text.\",\n# \"End of demonstration.\" and the trailing empty paragraph. Paragraph 3\n# keeps being the next one to remove as the collection re-indexes after\n# each delete.\n$d.Paragraphs.Item(3).Range.Delete()\n$d.Paragraphs.Item(3).Range.Delete()\n$d.Paragraphs.Item(3).Range.Delete()\n"}
